$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "70.778.77"
$ws.Range("E2").Value = "  +2.87%  "
$ws.Range("D3").Value = "3.792.43"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "701.82"
$ws.Range("E5").Value = "  +10.40%  "
$ws.Range("D6").Value = "173.09"
$ws.Range("E6").Value = "  +4.72%  "
$ws.Range("D7").Value = "3.790.84"
$ws.Range("E7").Value = "  +0.88%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.526"
$ws.Range("E9").Value = "  +1.00%  "
$ws.Range("D10").Value = "0.162"
$ws.Range("E10").Value = "  +2.68%  "
$ws.Range("D11").Value = "7.41"
$ws.Range("E11").Value = "  +7.36%  "
$ws.Range("D12").Value = "0.459"
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("E13").Value = "  +7.95%  "
$ws.Range("D14").Value = "36.30"
$ws.Range("E14").Value = "  +4.28%  "
$ws.Range("D15").Value = "4.434.76"
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").Value = "3.795.23"
$ws.Range("E16").Value = "  +0.94%  "
$ws.Range("D17").Value = "70.812.08"
$ws.Range("E17").Value = "  +2.90%  "
$ws.Range("D18").Value = "17.81"
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("D19").Value = "7.18"
$ws.Range("E19").Value = "  +2.96%  "
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("D21").Value = "11.10"
$ws.Range("E21").Value = "  +16.97%  "
$ws.Range("D22").Value = "481.37"
$ws.Range("E22").Value = "  +2.57%  "
$ws.Range("E23").Value = "  +1.75%  "
$ws.Range("D24").Value = "84.32"
$ws.Range("E24").Value = "  +3.47%  "
$ws.Range("D25").Value = "0.0000143"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("D26").Value = "12.38"
$ws.Range("E26").Value = "  +2.10%  "
$ws.Range("D27").Value = "2.17"
$ws.Range("E27").Value = "  +3.55%  "
$ws.Range("D28").Value = "10.43"
$ws.Range("E28").Value = "  +4.09%  "
$ws.Range("D29").Value = "3.945.33"
$ws.Range("E29").Value = "  +1.00%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").Value = "3.14"
$ws.Range("E31").Value = "  +17.11%  "
$ws.Range("D32").Value = "7.52"
$ws.Range("E32").Value = "  +6.05%  "
$ws.Range("D33").Value = "2.27"
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("D34").Value = "29.53"
$ws.Range("E34").Value = "  +4.06%  "
$ws.Range("D35").Value = "0.182"
$ws.Range("E35").Value = "  +4.79%  "
$ws.Range("D36").Value = "9.22"
$ws.Range("E36").Value = "  +4.17%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  +1.98%  "
$ws.Range("D39").Value = "3.43"
$ws.Range("E39").Value = "  +6.31%  "
$ws.Range("D40").Value = "6.03"
$ws.Range("E40").Value = "  +4.85%  "
$ws.Range("E41").Value = "  +12.69%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "0.974"
$ws.Range("E42").Value = "  +2.16%  "
$ws.Range("B43").Value = "FLOKI"
$ws.Range("C43").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D43").Value = "0.000325"
$ws.Range("E43").Value = "  +22.05%  "
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D46").Value = "162.37"
$ws.Range("E46").Value = "  +4.49%  "
$ws.Range("D47").Value = "49.00"
$ws.Range("E47").Value = "  +3.41%  "
$ws.Range("D48").Value = "44.74"
$ws.Range("E48").Value = "  +0.72%  "
$ws.Range("B49").Value = "TheGraph"
$ws.Range("C49").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D49").Value = "0.302"
$ws.Range("E49").Value = "  +3.44%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Value = "1.38"
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("D51").Value = "8.55"
$ws.Range("E51").Value = "  +2.50%  "
